# Updates cryptos list values (Price and Volume(1h) columns) to match
# the latest scraped data, as committed by the GitHub Actions workflow.
#
# Numeric-looking Price strings (e.g. "570.71") are written via a
# temporary Text number format so Excel's COM layer stores them as
# literal text (matching the source data) instead of silently coercing
# them to floating-point numbers; the style is reset back to Normal
# immediately after so no stray cell formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '62.672.34'
$ws.Range('E2').Value = '  -0.63%  '
$ws.Range('D3').Value = '2.454.56'
$ws.Range('E3').Value = '  -0.67%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.71'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.17'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.39%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.528'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.88%  '
$ws.Range('E9').Value = '  -1.32%  '
$ws.Range('E10').Value = '  -0.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.16'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.346'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.56'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.98%  '
$ws.Range('E14').Value = '  -3.47%  '
$ws.Range('D15').Value = '2.901.00'
$ws.Range('E15').Value = '  -0.66%  '
$ws.Range('D16').Value = '62.469.45'
$ws.Range('E16').Value = '  -0.87%  '
$ws.Range('D17').Value = '2.456.35'
$ws.Range('E17').Value = '  -0.69%  '
$ws.Range('E18').Value = '  -6.61%  '
$ws.Range('E19').Value = '  -3.18%  '
$ws.Range('E20').Value = '  -0.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '321.23'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.48%  '
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.90'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.73%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '64.68'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '645.56'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.66%  '
$ws.Range('D27').Value = '2.577.58'
$ws.Range('E27').Value = '  -0.60%  '
$ws.Range('E28').Value = '  -4.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('E30').Value = '  -3.41%  '
$ws.Range('E31').Value = '  -2.76%  '
$ws.Range('E32').Value = '  -3.32%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.132'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.77%  '
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('E35').Value = '  -4.38%  '
$ws.Range('E36').Value = '  -3.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '150.67'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.52'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.41%  '
$ws.Range('E39').Value = '  -2.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.30'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.82%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.63'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.91%  '
$ws.Range('E42').Value = '  -3.71%  '
$ws.Range('E43').Value = '  +2.42%  '
$ws.Range('E44').Value = '  +0.66%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '152.49'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.41%  '
$ws.Range('E46').Value = '  +1.73%  '
$ws.Range('E47').Value = '  -2.06%  '
$ws.Range('E48').Value = '  -0.60%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.93'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.31%  '
$ws.Range('E50').Value = '  -1.58%  '
$ws.Range('E51').Value = '  -1.85%  '
